# Applies "fixed integrations test refuel" commit to the "test" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# ---------------------------------------------------------------------
# 1. Relabel the "Aktuelt" column header to "Resultat" in every test
#    block except the "Test Gear_Up" one at row 40 (left untouched).
# ---------------------------------------------------------------------
foreach ($r in 16, 24, 32, 48, 56, 64) {
    $ws.Cells.Item($r, 5).Value = "Resultat"
}

# ---------------------------------------------------------------------
# 2. Rename the "Test Fuel_Empty Throttle_100%" block to
#    "Test Refuel Throttle_100%".
# ---------------------------------------------------------------------
$ws.Cells.Item(64, 1).Value = "Test Refuel Throttle_100%"

# ---------------------------------------------------------------------
# 3. Update the refuel test's recorded state (rows 66-70).
# ---------------------------------------------------------------------

# Row 66 - currentGear
$ws.Cells.Item(66, 2).ClearContents()       # B66 (Tilstand) now blank
$ws.Cells.Item(66, 4).Value = 2             # D66 (Forventet)
$ws.Cells.Item(66, 5).Value = 2             # E66 (Resultat)

# Row 67 - clutchEngaged
$ws.Cells.Item(67, 2).ClearContents()       # B67 (Tilstand) now blank

# Row 68 - throttlePosition
$ws.Cells.Item(68, 3).Value = 1             # C68 (Input)
$ws.Cells.Item(68, 4).Value = 0             # D68 (Forventet)
$ws.Cells.Item(68, 5).Value = 1             # E68 (Resultat)

# Row 69 - currentRpm
$ws.Cells.Item(69, 4).Value = 100           # D69 (Forventet)
$ws.Cells.Item(69, 5).Value = 100           # E69 (Resultat)

# Row 70 - f_content
$ws.Cells.Item(70, 2).ClearContents()       # B70 (Tilstand) now blank
$ws.Cells.Item(70, 3).Value = 100           # C70 (Input)

# D70 needs to become the text "99.75" (like D54) while keeping the
# same cell style as D54/D68 (s="11"), instead of its previous
# "#,##0" numeric style. Copy format+value from D54, which already
# holds that exact text, then restore D70's own address.
$ws.Range("D54").Copy() | Out-Null
$ws.Range("D70").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D54").Copy() | Out-Null
$ws.Range("D70").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item(70, 5).Value = 99.75         # E70 (Resultat)

# ---------------------------------------------------------------------
# 4. Drop the leftover blank rows 72-78 that trailed the sheet.
# ---------------------------------------------------------------------
$ws.Range("A72:E78").Clear()

# ---------------------------------------------------------------------
# 5. Restore the view state (selection / scroll position).
# ---------------------------------------------------------------------
$null = $ws.Range("B74").Select()
$excel.ActiveWindow.ScrollRow = 44
